# Auto-generated edit script: updates cryptocurrency price/name data in Sheet1
# to match the "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'282.24"
$ws.Range("D3").Value = "'20.54"
$ws.Range("D4").Value = "'6.248"
$ws.Range("D5").Value = "'0.06152"
$ws.Range("D7").Value = "'6.565"
$ws.Range("D8").Value = "'1.507"
$ws.Range("D9").Value = "'0.8201"
$ws.Range("D10").Value = "'0.01380"
$ws.Range("D11").Value = "'0.1637"
$ws.Range("D12").Value = "'0.08385"
$ws.Range("D13").Value = "'0.03531"
$ws.Range("D14").Value = "'0.03183"
$ws.Range("D15").Value = "'0.09129"
$ws.Range("D16").Value = "'3.704"
$ws.Range("D17").Value = "'0.001643"
$ws.Range("D18").Value = "'0.04702"
$ws.Range("D19").Value = "'0.006409"
$ws.Range("D20").Value = "'0.006154"
$ws.Range("D21").Value = "'0.001070"
$ws.Range("D22").Value = "'0.0001612"
$ws.Range("D23").Value = "'3.767"
$ws.Range("D24").Value = "'2.323"
$ws.Range("D25").Value = "'0.3356"
$ws.Range("D40").Value = "'0.04674"
$ws.Range("D41").Value = "'0.007206"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.004505"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1096"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01103"
$ws.Range("D45").Value = "'0.00006605"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'1.001"
$ws.Range("D48").Value = "'0.002948"
$ws.Range("D49").Value = "'0.00001902"
$ws.Range("D50").Value = "'0.01241"
